$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores values as text (e.g. "87.451.17", trailing
# zeros like "204.50", thousands-style dot grouping). Force the cells we are
# about to rewrite to Text format first so Excel does not auto-convert the
# numeric-looking strings into real numbers (which would drop trailing zeros,
# switch to scientific notation, etc).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D30", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "87.451.17"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "3.151.34"
$ws.Range("E3").Value = "  -6.25%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "204.50"
$ws.Range("E5").Value = "  -6.85%  "
$ws.Range("D6").Value = "606.58"
$ws.Range("E6").Value = "  -6.22%  "
$ws.Range("D7").Value = "0.378"
$ws.Range("E7").Value = "  -9.15%  "
$ws.Range("D8").Value = "0.659"
$ws.Range("E8").Value = "  +7.43%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "3.146.51"
$ws.Range("E10").Value = "  -6.38%  "
$ws.Range("D11").Value = "0.532"
$ws.Range("E11").Value = "  -13.55%  "
$ws.Range("D12").Value = "0.177"
$ws.Range("E12").Value = "  +5.18%  "
$ws.Range("D13").Value = "0.0000243"
$ws.Range("E13").Value = "  -16.75%  "
$ws.Range("D14").Value = "3.722.89"
$ws.Range("E14").Value = "  -6.42%  "
$ws.Range("D15").Value = "5.22"
$ws.Range("E15").Value = "  -5.75%  "
$ws.Range("D16").Value = "87.111.52"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "31.85"
$ws.Range("E17").Value = "  -12.82%  "
$ws.Range("D18").Value = "3.149.41"
$ws.Range("E18").Value = "  -6.42%  "
$ws.Range("D19").Value = "2.99"
$ws.Range("E19").Value = "  -4.57%  "
$ws.Range("D20").Value = "13.34"
$ws.Range("E20").Value = "  -9.98%  "
$ws.Range("D21").Value = "412.69"
$ws.Range("E21").Value = "  -9.99%  "
$ws.Range("D22").Value = "8.45"
$ws.Range("E22").Value = "  -12.96%  "
$ws.Range("D23").Value = "5.06"
$ws.Range("E23").Value = "  -8.88%  "
$ws.Range("D24").Value = "5.14"
$ws.Range("E24").Value = "  -7.52%  "
$ws.Range("D25").Value = "11.81"
$ws.Range("E25").Value = "  -7.59%  "
$ws.Range("D26").Value = "3.311.83"
$ws.Range("E26").Value = "  -6.03%  "
$ws.Range("D27").Value = "73.26"
$ws.Range("D28").Value = "0.0000129"
$ws.Range("E28").Value = "  -9.90%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").Value = "0.161"
$ws.Range("E30").Value = "  -19.45%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").Value = "534.99"
$ws.Range("E32").Value = "  -10.12%  "
$ws.Range("D33").Value = "8.24"
$ws.Range("E33").Value = "  -12.67%  "
$ws.Range("D34").Value = "1.30"
$ws.Range("E34").Value = "  -17.51%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value = "6.65"
$ws.Range("E35").Value = "  -9.08%  "
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "1.83"
$ws.Range("E36").Value = "  -13.18%  "
$ws.Range("E37").Value = "  -8.25%  "
$ws.Range("D38").Value = "21.82"
$ws.Range("E38").Value = "  -7.36%  "
$ws.Range("D39").Value = "21.79"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "2.98"
$ws.Range("E41").Value = "  -6.68%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").Value = "1.88"
$ws.Range("E43").Value = "  -12.70%  "
$ws.Range("D44").Value = "0.369"
$ws.Range("E44").Value = "  -14.07%  "
$ws.Range("D45").Value = "147.39"
$ws.Range("E45").Value = "  -6.79%  "
$ws.Range("D46").Value = "171.50"
$ws.Range("E46").Value = "  -9.46%  "
$ws.Range("D47").Value = "43.13"
$ws.Range("E47").Value = "  -7.12%  "
$ws.Range("D48").Value = "0.127"
$ws.Range("E48").Value = "  +5.60%  "
$ws.Range("E49").Value = "  -15.16%  "
$ws.Range("D50").Value = "3.95"
$ws.Range("E50").Value = "  -12.35%  "
$ws.Range("D51").Value = "0.691"
$ws.Range("E51").Value = "  -12.17%  "
